# Apply edit: add 'Ontogeny'/'Protein' data for TestPopulation and add a new
# population row 'TestPopulation_noOnto' (without ontogeny info) on the
# Demographics sheet (sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Fill in the Protein (Q) and Ontogeny (R) columns for the existing
# TestPopulation row (row 2).
$ws.Range("Q2").Value = "CYP3A4, CYP2D6"
$ws.Range("R2").Value = "CYP3A4, CYP2D6"

# Add a new row (row 3) describing a population without ontogeny info,
# with the same values as row 2 (columns A:P) except for the population
# name itself.
$ws.Range("A3").Value = "TestPopulation_noOnto"
$ws.Range("B3").Value = "Human"
$ws.Range("C3").Value = "European_ICRP_2002"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0
$ws.Range("H3").Value = "kg"
$ws.Range("K3").Value = "cm"
$ws.Range("L3").Value = 22
$ws.Range("M3").Value = 41
$ws.Range("P3").Value = "kg/m²"

$ws.Range("R3").Select()

$wb.Save()
